# Applies the "AI Energy Consumption Dataset" edit:
#  - Adds the efficiency-gains exponent as a standalone input cell (J20 = 0.9)
#  - Updates the label in I20 to drop the now-stale "(0.7)" hint
#  - Re-points every POWER(...) formula in columns H:M (rows 2,4,6,8,10,12)
#    to use that J20 input instead of the hard-coded 0.7 literal
#  - Moves the active selection to I4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New exponent input cell, used by all the POWER() formulas below.
$ws.Range("J20").Value = 0.9

# The "Formula Used :" helper label no longer needs to echo the exponent
# value inline, since it now lives in J20.
$ws.Range("I20").Value = "b - efficiency gains"

# Columns and their numerator constants, matching the original formulas:
#   H -> 1750 (ChatGPT-4), I -> 170 (ChatGPT-3.5), J -> 1500 (Gemini),
#   K -> 200 (Claude 3), L -> 70 (LLaMA 3), M -> 13 (LLaMA 2)
$numerators = @{ "H" = 1750; "I" = 170; "J" = 1500; "K" = 200; "L" = 70; "M" = 13 }
$rows = @(2, 4, 6, 8, 10, 12)

foreach ($r in $rows) {
    foreach ($col in @("H", "I", "J", "K", "L", "M")) {
        # Rows 2,4,6,8,10 have all six columns populated; row 12 only has
        # H and I as formulas (J12:M12 remain the literal "NIL" text).
        if ($r -eq 12 -and ($col -eq "J" -or $col -eq "K" -or $col -eq "L" -or $col -eq "M")) {
            continue
        }
        $num = $numerators[$col]
        $ws.Range("$col$r").Formula = "=(POWER(($num/G$r),J20)*E$r)"
    }
}

# Move the active selection, as recorded in the saved sheet view.
$ws.Range("I4").Select()
